# edit.ps1
# Applies the "script update" commit to the premier-league 2023-2024 match-odds sheet:
#   1) Several groups of adjacent rows had their match data (columns F:V) swapped /
#      re-ordered (the A:E "meta" columns - Indice/pais/torneio/temporada/data_partida -
#      stay untouched since the matches within a group share the same kickoff date).
#   2) Four brand-new match rows (117-120) are appended at the bottom, extending the
#      used range from A1:V116 to A1:V120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order swapped match rows (columns F:V) per the source permutation ---
# Row 4 <- old row 5
$ws.Cells.Item(4, 6).Value = "Brighton"
$ws.Cells.Item(4, 7).Value = 4
$ws.Cells.Item(4, 8).Value = "Luton"
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1.33
$ws.Cells.Item(4, 11).Value = "15/06/2023 10:29"
$ws.Cells.Item(4, 12).Value = 1.27
$ws.Cells.Item(4, 13).Value = "12/08/2023 15:37"
$ws.Cells.Item(4, 14).Value = 5.16
$ws.Cells.Item(4, 15).Value = "15/06/2023 10:29"
$ws.Cells.Item(4, 16).Value = 6.36
$ws.Cells.Item(4, 17).Value = "12/08/2023 15:57"
$ws.Cells.Item(4, 18).Value = 8.29
$ws.Cells.Item(4, 19).Value = "15/06/2023 10:29"
$ws.Cells.Item(4, 20).Value = 11.36
$ws.Cells.Item(4, 21).Value = "12/08/2023 15:57"
$ws.Cells.Item(4, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-luton/Sd9uKdMe/"

# Row 5 <- old row 6
$ws.Cells.Item(5, 6).Value = "Bournemouth"
$ws.Cells.Item(5, 7).Value = 1
$ws.Cells.Item(5, 8).Value = "West Ham"
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 2.64
$ws.Cells.Item(5, 11).Value = "15/06/2023 10:16"
$ws.Cells.Item(5, 12).Value = 2.75
$ws.Cells.Item(5, 13).Value = "12/08/2023 15:59"
$ws.Cells.Item(5, 14).Value = 3.38
$ws.Cells.Item(5, 15).Value = "15/06/2023 10:16"
$ws.Cells.Item(5, 16).Value = 3.6
$ws.Cells.Item(5, 17).Value = "12/08/2023 15:57"
$ws.Cells.Item(5, 18).Value = 2.64
$ws.Cells.Item(5, 19).Value = "15/06/2023 10:16"
$ws.Cells.Item(5, 20).Value = 2.63
$ws.Cells.Item(5, 21).Value = "12/08/2023 15:57"
$ws.Cells.Item(5, 22).Value = "https://www.betexplorer.com/football/england/premier-league/bournemouth-west-ham/YZ9yLx7k/"

# Row 6 <- old row 4
$ws.Cells.Item(6, 6).Value = "Sheffield Utd"
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = "Crystal Palace"
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 2.72
$ws.Cells.Item(6, 11).Value = "15/06/2023 10:28"
$ws.Cells.Item(6, 12).Value = 3.98
$ws.Cells.Item(6, 13).Value = "12/08/2023 15:59"
$ws.Cells.Item(6, 14).Value = 3.14
$ws.Cells.Item(6, 15).Value = "15/06/2023 10:28"
$ws.Cells.Item(6, 16).Value = 3.22
$ws.Cells.Item(6, 17).Value = "12/08/2023 15:54"
$ws.Cells.Item(6, 18).Value = 2.72
$ws.Cells.Item(6, 19).Value = "15/06/2023 10:28"
$ws.Cells.Item(6, 20).Value = 2.17
$ws.Cells.Item(6, 21).Value = "12/08/2023 15:54"
$ws.Cells.Item(6, 22).Value = "https://www.betexplorer.com/football/england/premier-league/sheffield-utd-crystal-palace/hjTJs0E2/"

# Row 23 <- old row 26
$ws.Cells.Item(23, 6).Value = "Manchester Utd"
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = "Nottingham"
$ws.Cells.Item(23, 9).Value = 2
$ws.Cells.Item(23, 10).Value = 1.26
$ws.Cells.Item(23, 11).Value = "13/08/2023 09:01"
$ws.Cells.Item(23, 12).Value = 1.38
$ws.Cells.Item(23, 13).Value = "26/08/2023 15:56"
$ws.Cells.Item(23, 14).Value = 6.29
$ws.Cells.Item(23, 15).Value = "13/08/2023 09:01"
$ws.Cells.Item(23, 16).Value = 5.3
$ws.Cells.Item(23, 17).Value = "26/08/2023 15:56"
$ws.Cells.Item(23, 18).Value = 11.33
$ws.Cells.Item(23, 19).Value = "13/08/2023 09:01"
$ws.Cells.Item(23, 20).Value = 8.51
$ws.Cells.Item(23, 21).Value = "26/08/2023 15:59"
$ws.Cells.Item(23, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-nottingham/4IjJ0EKs/"

# Row 24 <- old row 23
$ws.Cells.Item(24, 6).Value = "Everton"
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = "Wolves"
$ws.Cells.Item(24, 9).Value = 1
$ws.Cells.Item(24, 10).Value = 2.13
$ws.Cells.Item(24, 11).Value = "13/08/2023 09:01"
$ws.Cells.Item(24, 12).Value = 2.45
$ws.Cells.Item(24, 13).Value = "26/08/2023 15:59"
$ws.Cells.Item(24, 14).Value = 3.32
$ws.Cells.Item(24, 15).Value = "13/08/2023 09:01"
$ws.Cells.Item(24, 16).Value = 3.38
$ws.Cells.Item(24, 17).Value = "26/08/2023 15:43"
$ws.Cells.Item(24, 18).Value = 3.78
$ws.Cells.Item(24, 19).Value = "13/08/2023 09:01"
$ws.Cells.Item(24, 20).Value = 3.14
$ws.Cells.Item(24, 21).Value = "26/08/2023 15:59"
$ws.Cells.Item(24, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-wolves/294i5fCQ/"

# Row 25 <- old row 24
$ws.Cells.Item(25, 6).Value = "Brentford"
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = "Crystal Palace"
$ws.Cells.Item(25, 9).Value = 1
$ws.Cells.Item(25, 10).Value = 1.93
$ws.Cells.Item(25, 11).Value = "13/08/2023 09:01"
$ws.Cells.Item(25, 12).Value = 2.16
$ws.Cells.Item(25, 13).Value = "26/08/2023 15:55"
$ws.Cells.Item(25, 14).Value = 3.49
$ws.Cells.Item(25, 15).Value = "13/08/2023 09:01"
$ws.Cells.Item(25, 16).Value = 3.37
$ws.Cells.Item(25, 17).Value = "26/08/2023 15:51"
$ws.Cells.Item(25, 18).Value = 4.3
$ws.Cells.Item(25, 19).Value = "13/08/2023 09:01"
$ws.Cells.Item(25, 20).Value = 3.81
$ws.Cells.Item(25, 21).Value = "26/08/2023 15:57"
$ws.Cells.Item(25, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-crystal-palace/EmOdqDk7/"

# Row 26 <- old row 25
$ws.Cells.Item(26, 6).Value = "Arsenal"
$ws.Cells.Item(26, 7).Value = 2
$ws.Cells.Item(26, 8).Value = "Fulham"
$ws.Cells.Item(26, 9).Value = 2
$ws.Cells.Item(26, 10).Value = 1.26
$ws.Cells.Item(26, 11).Value = "13/08/2023 09:01"
$ws.Cells.Item(26, 12).Value = 1.29
$ws.Cells.Item(26, 13).Value = "26/08/2023 15:39"
$ws.Cells.Item(26, 14).Value = 6.29
$ws.Cells.Item(26, 15).Value = "13/08/2023 09:01"
$ws.Cells.Item(26, 16).Value = 6.29
$ws.Cells.Item(26, 17).Value = "26/08/2023 15:44"
$ws.Cells.Item(26, 18).Value = 11.33
$ws.Cells.Item(26, 19).Value = "13/08/2023 09:01"
$ws.Cells.Item(26, 20).Value = 10.32
$ws.Cells.Item(26, 21).Value = "26/08/2023 15:44"
$ws.Cells.Item(26, 22).Value = "https://www.betexplorer.com/football/england/premier-league/arsenal-fulham/6NOhpgz1/"

# Row 33 <- old row 35
$ws.Cells.Item(33, 6).Value = "Chelsea"
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = "Nottingham"
$ws.Cells.Item(33, 9).Value = 1
$ws.Cells.Item(33, 10).Value = 1.51
$ws.Cells.Item(33, 11).Value = "20/08/2023 09:02"
$ws.Cells.Item(33, 12).Value = 1.41
$ws.Cells.Item(33, 13).Value = "02/09/2023 15:50"
$ws.Cells.Item(33, 14).Value = 4.47
$ws.Cells.Item(33, 15).Value = "20/08/2023 09:02"
$ws.Cells.Item(33, 16).Value = 5.01
$ws.Cells.Item(33, 17).Value = "02/09/2023 15:58"
$ws.Cells.Item(33, 18).Value = 6.71
$ws.Cells.Item(33, 19).Value = "20/08/2023 09:02"
$ws.Cells.Item(33, 20).Value = 8.45
$ws.Cells.Item(33, 21).Value = "02/09/2023 15:58"
$ws.Cells.Item(33, 22).Value = "https://www.betexplorer.com/football/england/premier-league/chelsea-nottingham/0d8k37tt/"

# Row 34 <- old row 36
$ws.Cells.Item(34, 6).Value = "Manchester City"
$ws.Cells.Item(34, 7).Value = 5
$ws.Cells.Item(34, 8).Value = "Fulham"
$ws.Cells.Item(34, 9).Value = 1
$ws.Cells.Item(34, 10).Value = 1.15
$ws.Cells.Item(34, 11).Value = "20/08/2023 09:02"
$ws.Cells.Item(34, 12).Value = 1.19
$ws.Cells.Item(34, 13).Value = "02/09/2023 15:51"
$ws.Cells.Item(34, 14).Value = 9.19
$ws.Cells.Item(34, 15).Value = "20/08/2023 09:02"
$ws.Cells.Item(34, 16).Value = 7.65
$ws.Cells.Item(34, 17).Value = "02/09/2023 15:54"
$ws.Cells.Item(34, 18).Value = 18
$ws.Cells.Item(34, 19).Value = "20/08/2023 09:02"
$ws.Cells.Item(34, 20).Value = 16
$ws.Cells.Item(34, 21).Value = "02/09/2023 15:58"
$ws.Cells.Item(34, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-city-fulham/Uy06aPu5/"

# Row 35 <- old row 33
$ws.Cells.Item(35, 6).Value = "Brentford"
$ws.Cells.Item(35, 7).Value = 2
$ws.Cells.Item(35, 8).Value = "Bournemouth"
$ws.Cells.Item(35, 9).Value = 2
$ws.Cells.Item(35, 10).Value = 1.73
$ws.Cells.Item(35, 11).Value = "20/08/2023 09:02"
$ws.Cells.Item(35, 12).Value = 1.75
$ws.Cells.Item(35, 13).Value = "02/09/2023 15:27"
$ws.Cells.Item(35, 14).Value = 4.37
$ws.Cells.Item(35, 15).Value = "20/08/2023 09:02"
$ws.Cells.Item(35, 16).Value = 4.01
$ws.Cells.Item(35, 17).Value = "02/09/2023 15:30"
$ws.Cells.Item(35, 18).Value = 3.95
$ws.Cells.Item(35, 19).Value = "20/08/2023 09:02"
$ws.Cells.Item(35, 20).Value = 4.86
$ws.Cells.Item(35, 21).Value = "02/09/2023 15:59"
$ws.Cells.Item(35, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-bournemouth/hh2ZdWJ6/"

# Row 36 <- old row 34
$ws.Cells.Item(36, 6).Value = "Burnley"
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(36, 8).Value = "Tottenham"
$ws.Cells.Item(36, 9).Value = 5
$ws.Cells.Item(36, 10).Value = 3.81
$ws.Cells.Item(36, 11).Value = "21/08/2023 06:14"
$ws.Cells.Item(36, 12).Value = 3.71
$ws.Cells.Item(36, 13).Value = "02/09/2023 15:53"
$ws.Cells.Item(36, 14).Value = 3.57
$ws.Cells.Item(36, 15).Value = "21/08/2023 06:14"
$ws.Cells.Item(36, 16).Value = 3.79
$ws.Cells.Item(36, 17).Value = "02/09/2023 15:59"
$ws.Cells.Item(36, 18).Value = 1.94
$ws.Cells.Item(36, 19).Value = "21/08/2023 06:14"
$ws.Cells.Item(36, 20).Value = 2.04
$ws.Cells.Item(36, 21).Value = "02/09/2023 15:52"
$ws.Cells.Item(36, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-tottenham/E7jreAlJ/"

# Row 42 <- old row 45
$ws.Cells.Item(42, 6).Value = "Tottenham"
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = "Sheffield Utd"
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = 1.29
$ws.Cells.Item(42, 11).Value = "28/08/2023 11:55"
$ws.Cells.Item(42, 12).Value = 1.3
$ws.Cells.Item(42, 13).Value = "16/09/2023 15:58"
$ws.Cells.Item(42, 14).Value = 5.56
$ws.Cells.Item(42, 15).Value = "28/08/2023 11:55"
$ws.Cells.Item(42, 16).Value = 6.31
$ws.Cells.Item(42, 17).Value = "16/09/2023 15:58"
$ws.Cells.Item(42, 18).Value = 9.48
$ws.Cells.Item(42, 19).Value = "28/08/2023 11:55"
$ws.Cells.Item(42, 20).Value = 9.4
$ws.Cells.Item(42, 21).Value = "16/09/2023 15:58"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/england/premier-league/tottenham-sheffield-utd/nZoxDrA4/"

# Row 43 <- old row 42
$ws.Cells.Item(43, 6).Value = "Aston Villa"
$ws.Cells.Item(43, 7).Value = 3
$ws.Cells.Item(43, 8).Value = "Crystal Palace"
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 1.88
$ws.Cells.Item(43, 11).Value = "28/08/2023 09:02"
$ws.Cells.Item(43, 12).Value = 1.98
$ws.Cells.Item(43, 13).Value = "16/09/2023 15:53"
$ws.Cells.Item(43, 14).Value = 3.72
$ws.Cells.Item(43, 15).Value = "28/08/2023 09:02"
$ws.Cells.Item(43, 16).Value = 3.66
$ws.Cells.Item(43, 17).Value = "16/09/2023 15:53"
$ws.Cells.Item(43, 18).Value = 3.91
$ws.Cells.Item(43, 19).Value = "28/08/2023 09:02"
$ws.Cells.Item(43, 20).Value = 4.08
$ws.Cells.Item(43, 21).Value = "16/09/2023 15:53"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/england/premier-league/aston-villa-crystal-palace/SUEBdNPN/"

# Row 44 <- old row 43
$ws.Cells.Item(44, 6).Value = "Fulham"
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = "Luton"
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 1.68
$ws.Cells.Item(44, 11).Value = "28/08/2023 11:54"
$ws.Cells.Item(44, 12).Value = 1.69
$ws.Cells.Item(44, 13).Value = "16/09/2023 15:36"
$ws.Cells.Item(44, 14).Value = 3.8
$ws.Cells.Item(44, 15).Value = "28/08/2023 11:54"
$ws.Cells.Item(44, 16).Value = 3.82
$ws.Cells.Item(44, 17).Value = "16/09/2023 15:59"
$ws.Cells.Item(44, 18).Value = 4.97
$ws.Cells.Item(44, 19).Value = "28/08/2023 11:54"
$ws.Cells.Item(44, 20).Value = 5.68
$ws.Cells.Item(44, 21).Value = "16/09/2023 15:59"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/england/premier-league/fulham-luton/bD5si1mo/"

# Row 45 <- old row 44
$ws.Cells.Item(45, 6).Value = "Manchester Utd"
$ws.Cells.Item(45, 7).Value = 1
$ws.Cells.Item(45, 8).Value = "Brighton"
$ws.Cells.Item(45, 9).Value = 3
$ws.Cells.Item(45, 10).Value = 1.88
$ws.Cells.Item(45, 11).Value = "28/08/2023 09:02"
$ws.Cells.Item(45, 12).Value = 2.07
$ws.Cells.Item(45, 13).Value = "16/09/2023 15:59"
$ws.Cells.Item(45, 14).Value = 3.99
$ws.Cells.Item(45, 15).Value = "28/08/2023 09:02"
$ws.Cells.Item(45, 16).Value = 3.9
$ws.Cells.Item(45, 17).Value = "16/09/2023 15:59"
$ws.Cells.Item(45, 18).Value = 3.6
$ws.Cells.Item(45, 19).Value = "28/08/2023 09:02"
$ws.Cells.Item(45, 20).Value = 3.53
$ws.Cells.Item(45, 21).Value = "16/09/2023 15:59"
$ws.Cells.Item(45, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-brighton/0IqQFpQo/"

# Row 51 <- old row 53
$ws.Cells.Item(51, 6).Value = "Luton"
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = "Wolves"
$ws.Cells.Item(51, 9).Value = 1
$ws.Cells.Item(51, 10).Value = 3.36
$ws.Cells.Item(51, 11).Value = "05/09/2023 12:11"
$ws.Cells.Item(51, 12).Value = 3.28
$ws.Cells.Item(51, 13).Value = "23/09/2023 15:59"
$ws.Cells.Item(51, 14).Value = 3.4
$ws.Cells.Item(51, 15).Value = "05/09/2023 12:11"
$ws.Cells.Item(51, 16).Value = 3.31
$ws.Cells.Item(51, 17).Value = "23/09/2023 15:59"
$ws.Cells.Item(51, 18).Value = 2.15
$ws.Cells.Item(51, 19).Value = "05/09/2023 12:11"
$ws.Cells.Item(51, 20).Value = 2.41
$ws.Cells.Item(51, 21).Value = "23/09/2023 15:59"
$ws.Cells.Item(51, 22).Value = "https://www.betexplorer.com/football/england/premier-league/luton-wolves/Mi0KPLgA/"

# Row 53 <- old row 51
$ws.Cells.Item(53, 6).Value = "Manchester City"
$ws.Cells.Item(53, 7).Value = 2
$ws.Cells.Item(53, 8).Value = "Nottingham"
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 1.12
$ws.Cells.Item(53, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(53, 12).Value = 1.17
$ws.Cells.Item(53, 13).Value = "23/09/2023 15:28"
$ws.Cells.Item(53, 14).Value = 9.08
$ws.Cells.Item(53, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(53, 16).Value = 8.5
$ws.Cells.Item(53, 17).Value = "23/09/2023 15:17"
$ws.Cells.Item(53, 18).Value = 15.94
$ws.Cells.Item(53, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(53, 20).Value = 18
$ws.Cells.Item(53, 21).Value = "23/09/2023 15:28"
$ws.Cells.Item(53, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-city-nottingham/GSENOu9G/"

# Row 56 <- old row 57
$ws.Cells.Item(56, 6).Value = "Brighton"
$ws.Cells.Item(56, 7).Value = 3
$ws.Cells.Item(56, 8).Value = "Bournemouth"
$ws.Cells.Item(56, 9).Value = 1
$ws.Cells.Item(56, 10).Value = 1.32
$ws.Cells.Item(56, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(56, 12).Value = 1.52
$ws.Cells.Item(56, 13).Value = "24/09/2023 14:30"
$ws.Cells.Item(56, 14).Value = 6.03
$ws.Cells.Item(56, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(56, 16).Value = 4.97
$ws.Cells.Item(56, 17).Value = "24/09/2023 14:53"
$ws.Cells.Item(56, 18).Value = 8.84
$ws.Cells.Item(56, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(56, 20).Value = 5.9
$ws.Cells.Item(56, 21).Value = "24/09/2023 14:59"
$ws.Cells.Item(56, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-bournemouth/fuLL4KHp/"

# Row 57 <- old row 56
$ws.Cells.Item(57, 6).Value = "Liverpool"
$ws.Cells.Item(57, 7).Value = 3
$ws.Cells.Item(57, 8).Value = "West Ham"
$ws.Cells.Item(57, 9).Value = 1
$ws.Cells.Item(57, 10).Value = 1.31
$ws.Cells.Item(57, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(57, 12).Value = 1.39
$ws.Cells.Item(57, 13).Value = "24/09/2023 14:53"
$ws.Cells.Item(57, 14).Value = 5.93
$ws.Cells.Item(57, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(57, 16).Value = 5.64
$ws.Cells.Item(57, 17).Value = "24/09/2023 14:58"
$ws.Cells.Item(57, 18).Value = 7.64
$ws.Cells.Item(57, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(57, 20).Value = 7.61
$ws.Cells.Item(57, 21).Value = "24/09/2023 14:59"
$ws.Cells.Item(57, 22).Value = "https://www.betexplorer.com/football/england/premier-league/liverpool-west-ham/r11GQ1v4/"

# Row 62 <- old row 63
$ws.Cells.Item(62, 6).Value = "Newcastle"
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = "Burnley"
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 1.34
$ws.Cells.Item(62, 11).Value = "23/09/2023 17:43"
$ws.Cells.Item(62, 12).Value = 1.38
$ws.Cells.Item(62, 13).Value = "30/09/2023 15:56"
$ws.Cells.Item(62, 14).Value = 5.32
$ws.Cells.Item(62, 15).Value = "23/09/2023 17:43"
$ws.Cells.Item(62, 16).Value = 5.24
$ws.Cells.Item(62, 17).Value = "30/09/2023 15:45"
$ws.Cells.Item(62, 18).Value = 7.68
$ws.Cells.Item(62, 19).Value = "23/09/2023 17:43"
$ws.Cells.Item(62, 20).Value = 8.53
$ws.Cells.Item(62, 21).Value = "30/09/2023 15:58"
$ws.Cells.Item(62, 22).Value = "https://www.betexplorer.com/football/england/premier-league/newcastle-utd-burnley/4hNjOTZ2/"

# Row 63 <- old row 64
$ws.Cells.Item(63, 6).Value = "Manchester Utd"
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = "Crystal Palace"
$ws.Cells.Item(63, 9).Value = 1
$ws.Cells.Item(63, 10).Value = 1.5
$ws.Cells.Item(63, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(63, 12).Value = 1.58
$ws.Cells.Item(63, 13).Value = "30/09/2023 15:58"
$ws.Cells.Item(63, 14).Value = 4.63
$ws.Cells.Item(63, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(63, 16).Value = 4.37
$ws.Cells.Item(63, 17).Value = "30/09/2023 15:59"
$ws.Cells.Item(63, 18).Value = 6.58
$ws.Cells.Item(63, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(63, 20).Value = 5.96
$ws.Cells.Item(63, 21).Value = "30/09/2023 15:59"
$ws.Cells.Item(63, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-crystal-palace/Q1OnP9Kd/"

# Row 64 <- old row 67
$ws.Cells.Item(64, 6).Value = "Everton"
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = "Luton"
$ws.Cells.Item(64, 9).Value = 2
$ws.Cells.Item(64, 10).Value = 1.58
$ws.Cells.Item(64, 11).Value = "23/09/2023 17:43"
$ws.Cells.Item(64, 12).Value = 1.67
$ws.Cells.Item(64, 13).Value = "30/09/2023 15:59"
$ws.Cells.Item(64, 14).Value = 3.94
$ws.Cells.Item(64, 15).Value = "23/09/2023 17:43"
$ws.Cells.Item(64, 16).Value = 3.93
$ws.Cells.Item(64, 17).Value = "30/09/2023 15:54"
$ws.Cells.Item(64, 18).Value = 5.63
$ws.Cells.Item(64, 19).Value = "23/09/2023 17:43"
$ws.Cells.Item(64, 20).Value = 5.7
$ws.Cells.Item(64, 21).Value = "30/09/2023 15:59"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-luton/8Qg2Hc1j/"

# Row 65 <- old row 62
$ws.Cells.Item(65, 6).Value = "West Ham"
$ws.Cells.Item(65, 7).Value = 2
$ws.Cells.Item(65, 8).Value = "Sheffield Utd"
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 1.49
$ws.Cells.Item(65, 11).Value = "23/09/2023 17:42"
$ws.Cells.Item(65, 12).Value = 1.48
$ws.Cells.Item(65, 13).Value = "30/09/2023 15:46"
$ws.Cells.Item(65, 14).Value = 4.43
$ws.Cells.Item(65, 15).Value = "23/09/2023 17:42"
$ws.Cells.Item(65, 16).Value = 4.79
$ws.Cells.Item(65, 17).Value = "30/09/2023 15:55"
$ws.Cells.Item(65, 18).Value = 5.97
$ws.Cells.Item(65, 19).Value = "23/09/2023 17:42"
$ws.Cells.Item(65, 20).Value = 7
$ws.Cells.Item(65, 21).Value = "30/09/2023 15:59"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/england/premier-league/west-ham-sheffield-utd/hEF3LRJL/"

# Row 66 <- old row 65
$ws.Cells.Item(66, 6).Value = "Wolves"
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = "Manchester City"
$ws.Cells.Item(66, 9).Value = 1
$ws.Cells.Item(66, 10).Value = 6.97
$ws.Cells.Item(66, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(66, 12).Value = 8.61
$ws.Cells.Item(66, 13).Value = "30/09/2023 15:58"
$ws.Cells.Item(66, 14).Value = 4.89
$ws.Cells.Item(66, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(66, 16).Value = 5.19
$ws.Cells.Item(66, 17).Value = "30/09/2023 15:58"
$ws.Cells.Item(66, 18).Value = 1.45
$ws.Cells.Item(66, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(66, 20).Value = 1.39
$ws.Cells.Item(66, 21).Value = "30/09/2023 15:58"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/england/premier-league/wolves-manchester-city/W4E7KoZR/"

# Row 67 <- old row 66
$ws.Cells.Item(67, 6).Value = "Bournemouth"
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = "Arsenal"
$ws.Cells.Item(67, 9).Value = 4
$ws.Cells.Item(67, 10).Value = 6.39
$ws.Cells.Item(67, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(67, 12).Value = 6.06
$ws.Cells.Item(67, 13).Value = "30/09/2023 15:58"
$ws.Cells.Item(67, 14).Value = 4.83
$ws.Cells.Item(67, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(67, 16).Value = 4.62
$ws.Cells.Item(67, 17).Value = "30/09/2023 15:59"
$ws.Cells.Item(67, 18).Value = 1.49
$ws.Cells.Item(67, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(67, 20).Value = 1.55
$ws.Cells.Item(67, 21).Value = "30/09/2023 15:53"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/england/premier-league/bournemouth-arsenal/xr3WMJwT/"

# Row 73 <- old row 74
$ws.Cells.Item(73, 6).Value = "Burnley"
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = "Chelsea"
$ws.Cells.Item(73, 9).Value = 4
$ws.Cells.Item(73, 10).Value = 4.24
$ws.Cells.Item(73, 11).Value = "28/09/2023 14:25"
$ws.Cells.Item(73, 12).Value = 5.2
$ws.Cells.Item(73, 13).Value = "07/10/2023 16:00"
$ws.Cells.Item(73, 14).Value = 3.78
$ws.Cells.Item(73, 15).Value = "28/09/2023 14:25"
$ws.Cells.Item(73, 16).Value = 4.07
$ws.Cells.Item(73, 17).Value = "07/10/2023 15:58"
$ws.Cells.Item(73, 18).Value = 1.79
$ws.Cells.Item(73, 19).Value = "28/09/2023 14:25"
$ws.Cells.Item(73, 20).Value = 1.69
$ws.Cells.Item(73, 21).Value = "07/10/2023 15:58"
$ws.Cells.Item(73, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-chelsea/pCfrEqCe/"

# Row 74 <- old row 75
$ws.Cells.Item(74, 6).Value = "Everton"
$ws.Cells.Item(74, 7).Value = 3
$ws.Cells.Item(74, 8).Value = "Bournemouth"
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 1.69
$ws.Cells.Item(74, 11).Value = "24/09/2023 10:02"
$ws.Cells.Item(74, 12).Value = 1.93
$ws.Cells.Item(74, 13).Value = "07/10/2023 15:58"
$ws.Cells.Item(74, 14).Value = 4.01
$ws.Cells.Item(74, 15).Value = "24/09/2023 10:02"
$ws.Cells.Item(74, 16).Value = 3.78
$ws.Cells.Item(74, 17).Value = "07/10/2023 15:58"
$ws.Cells.Item(74, 18).Value = 4.57
$ws.Cells.Item(74, 19).Value = "24/09/2023 10:02"
$ws.Cells.Item(74, 20).Value = 4.18
$ws.Cells.Item(74, 21).Value = "07/10/2023 15:58"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-bournemouth/CInUym42/"

# Row 75 <- old row 76
$ws.Cells.Item(75, 6).Value = "Manchester Utd"
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = "Brentford"
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 1.55
$ws.Cells.Item(75, 11).Value = "24/09/2023 10:01"
$ws.Cells.Item(75, 12).Value = 1.63
$ws.Cells.Item(75, 13).Value = "07/10/2023 15:50"
$ws.Cells.Item(75, 14).Value = 4.44
$ws.Cells.Item(75, 15).Value = "24/09/2023 10:01"
$ws.Cells.Item(75, 16).Value = 4.33
$ws.Cells.Item(75, 17).Value = "07/10/2023 15:53"
$ws.Cells.Item(75, 18).Value = 5.18
$ws.Cells.Item(75, 19).Value = "24/09/2023 10:01"
$ws.Cells.Item(75, 20).Value = 5.41
$ws.Cells.Item(75, 21).Value = "07/10/2023 15:54"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-brentford/8pQbIb3s/"

# Row 76 <- old row 73
$ws.Cells.Item(76, 6).Value = "Fulham"
$ws.Cells.Item(76, 7).Value = 3
$ws.Cells.Item(76, 8).Value = "Sheffield Utd"
$ws.Cells.Item(76, 9).Value = 1
$ws.Cells.Item(76, 10).Value = 1.61
$ws.Cells.Item(76, 11).Value = "28/09/2023 14:26"
$ws.Cells.Item(76, 12).Value = 1.52
$ws.Cells.Item(76, 13).Value = "07/10/2023 15:41"
$ws.Cells.Item(76, 14).Value = 3.93
$ws.Cells.Item(76, 15).Value = "28/09/2023 14:26"
$ws.Cells.Item(76, 16).Value = 4.56
$ws.Cells.Item(76, 17).Value = "07/10/2023 15:58"
$ws.Cells.Item(76, 18).Value = 5.39
$ws.Cells.Item(76, 19).Value = "28/09/2023 14:26"
$ws.Cells.Item(76, 20).Value = 6.57
$ws.Cells.Item(76, 21).Value = "07/10/2023 15:58"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/england/premier-league/fulham-sheffield-utd/j9oYz7J8/"

# Row 78 <- old row 80
$ws.Cells.Item(78, 6).Value = "Wolves"
$ws.Cells.Item(78, 7).Value = 1
$ws.Cells.Item(78, 8).Value = "Aston Villa"
$ws.Cells.Item(78, 9).Value = 1
$ws.Cells.Item(78, 10).Value = 2.81
$ws.Cells.Item(78, 11).Value = "24/09/2023 10:02"
$ws.Cells.Item(78, 12).Value = 3.7
$ws.Cells.Item(78, 13).Value = "08/10/2023 14:45"
$ws.Cells.Item(78, 14).Value = 3.39
$ws.Cells.Item(78, 15).Value = "24/09/2023 10:02"
$ws.Cells.Item(78, 16).Value = 3.68
$ws.Cells.Item(78, 17).Value = "08/10/2023 14:45"
$ws.Cells.Item(78, 18).Value = 2.62
$ws.Cells.Item(78, 19).Value = "24/09/2023 10:02"
$ws.Cells.Item(78, 20).Value = 2.08
$ws.Cells.Item(78, 21).Value = "08/10/2023 14:45"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/england/premier-league/wolves-aston-villa/GAT6GxYg/"

# Row 80 <- old row 78
$ws.Cells.Item(80, 6).Value = "Brighton"
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = "Liverpool"
$ws.Cells.Item(80, 9).Value = 2
$ws.Cells.Item(80, 10).Value = 2.47
$ws.Cells.Item(80, 11).Value = "24/09/2023 10:01"
$ws.Cells.Item(80, 12).Value = 3.26
$ws.Cells.Item(80, 13).Value = "08/10/2023 14:44"
$ws.Cells.Item(80, 14).Value = 3.9
$ws.Cells.Item(80, 15).Value = "24/09/2023 10:01"
$ws.Cells.Item(80, 16).Value = 4.19
$ws.Cells.Item(80, 17).Value = "08/10/2023 14:44"
$ws.Cells.Item(80, 18).Value = 2.54
$ws.Cells.Item(80, 19).Value = "24/09/2023 10:01"
$ws.Cells.Item(80, 20).Value = 2.1
$ws.Cells.Item(80, 21).Value = "08/10/2023 14:44"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-liverpool/2m5wFPdk/"

# Row 103 <- old row 105
$ws.Cells.Item(103, 6).Value = "Brentford"
$ws.Cells.Item(103, 7).Value = 3
$ws.Cells.Item(103, 8).Value = "West Ham"
$ws.Cells.Item(103, 9).Value = 2
$ws.Cells.Item(103, 10).Value = 1.95
$ws.Cells.Item(103, 11).Value = "21/10/2023 20:02"
$ws.Cells.Item(103, 12).Value = 2.16
$ws.Cells.Item(103, 13).Value = "04/11/2023 15:50"
$ws.Cells.Item(103, 14).Value = 3.65
$ws.Cells.Item(103, 15).Value = "21/10/2023 20:02"
$ws.Cells.Item(103, 16).Value = 3.71
$ws.Cells.Item(103, 17).Value = "04/11/2023 15:50"
$ws.Cells.Item(103, 18).Value = 4
$ws.Cells.Item(103, 19).Value = "21/10/2023 20:02"
$ws.Cells.Item(103, 20).Value = 3.39
$ws.Cells.Item(103, 21).Value = "04/11/2023 15:50"
$ws.Cells.Item(103, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-west-ham/MkBzuDGB/"

# Row 104 <- old row 106
$ws.Cells.Item(104, 6).Value = "Burnley"
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = "Crystal Palace"
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = 2.77
$ws.Cells.Item(104, 11).Value = "23/10/2023 15:48"
$ws.Cells.Item(104, 12).Value = 3.53
$ws.Cells.Item(104, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(104, 14).Value = 3.22
$ws.Cells.Item(104, 15).Value = "23/10/2023 15:48"
$ws.Cells.Item(104, 16).Value = 3.21
$ws.Cells.Item(104, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(104, 18).Value = 2.61
$ws.Cells.Item(104, 19).Value = "23/10/2023 15:48"
$ws.Cells.Item(104, 20).Value = 2.31
$ws.Cells.Item(104, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(104, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-crystal-palace/0tAvvXVH/"

# Row 105 <- old row 107
$ws.Cells.Item(105, 6).Value = "Everton"
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = "Brighton"
$ws.Cells.Item(105, 9).Value = 1
$ws.Cells.Item(105, 10).Value = 2.86
$ws.Cells.Item(105, 11).Value = "21/10/2023 20:02"
$ws.Cells.Item(105, 12).Value = 2.86
$ws.Cells.Item(105, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(105, 14).Value = 3.84
$ws.Cells.Item(105, 15).Value = "21/10/2023 20:02"
$ws.Cells.Item(105, 16).Value = 3.52
$ws.Cells.Item(105, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(105, 18).Value = 2.24
$ws.Cells.Item(105, 19).Value = "21/10/2023 20:02"
$ws.Cells.Item(105, 20).Value = 2.54
$ws.Cells.Item(105, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-brighton/f39rwioO/"

# Row 106 <- old row 104
$ws.Cells.Item(106, 6).Value = "Manchester City"
$ws.Cells.Item(106, 7).Value = 6
$ws.Cells.Item(106, 8).Value = "Bournemouth"
$ws.Cells.Item(106, 9).Value = 1
$ws.Cells.Item(106, 10).Value = 1.15
$ws.Cells.Item(106, 11).Value = "21/10/2023 22:01"
$ws.Cells.Item(106, 12).Value = 1.09
$ws.Cells.Item(106, 13).Value = "04/11/2023 15:50"
$ws.Cells.Item(106, 14).Value = 9.18
$ws.Cells.Item(106, 15).Value = "21/10/2023 22:01"
$ws.Cells.Item(106, 16).Value = 11.5
$ws.Cells.Item(106, 17).Value = "04/11/2023 15:21"
$ws.Cells.Item(106, 18).Value = 17.1
$ws.Cells.Item(106, 19).Value = "21/10/2023 22:01"
$ws.Cells.Item(106, 20).Value = 28.5
$ws.Cells.Item(106, 21).Value = "04/11/2023 15:52"
$ws.Cells.Item(106, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-city-bournemouth/AiwcdEon/"

# Row 107 <- old row 103
$ws.Cells.Item(107, 6).Value = "Sheffield Utd"
$ws.Cells.Item(107, 7).Value = 2
$ws.Cells.Item(107, 8).Value = "Wolves"
$ws.Cells.Item(107, 9).Value = 1
$ws.Cells.Item(107, 10).Value = 3.2
$ws.Cells.Item(107, 11).Value = "23/10/2023 15:48"
$ws.Cells.Item(107, 12).Value = 4.35
$ws.Cells.Item(107, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(107, 14).Value = 3.43
$ws.Cells.Item(107, 15).Value = "23/10/2023 15:48"
$ws.Cells.Item(107, 16).Value = 3.76
$ws.Cells.Item(107, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(107, 18).Value = 2.22
$ws.Cells.Item(107, 19).Value = "23/10/2023 15:48"
$ws.Cells.Item(107, 20).Value = 1.87
$ws.Cells.Item(107, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(107, 22).Value = "https://www.betexplorer.com/football/england/premier-league/sheffield-utd-wolves/0tW9gCV4/"

# Row 113 <- old row 114
$ws.Cells.Item(113, 6).Value = "Arsenal"
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = "Burnley"
$ws.Cells.Item(113, 9).Value = 1
$ws.Cells.Item(113, 10).Value = 1.24
$ws.Cells.Item(113, 11).Value = "29/10/2023 11:22"
$ws.Cells.Item(113, 12).Value = 1.19
$ws.Cells.Item(113, 13).Value = "11/11/2023 15:13"
$ws.Cells.Item(113, 14).Value = 6.5
$ws.Cells.Item(113, 15).Value = "29/10/2023 11:22"
$ws.Cells.Item(113, 16).Value = 7.49
$ws.Cells.Item(113, 17).Value = "11/11/2023 15:36"
$ws.Cells.Item(113, 18).Value = 12.06
$ws.Cells.Item(113, 19).Value = "29/10/2023 11:22"
$ws.Cells.Item(113, 20).Value = 16.59
$ws.Cells.Item(113, 21).Value = "11/11/2023 15:36"
$ws.Cells.Item(113, 22).Value = "https://www.betexplorer.com/football/england/premier-league/arsenal-burnley/ncYLjAFN/"

# Row 114 <- old row 113
$ws.Cells.Item(114, 6).Value = "Crystal Palace"
$ws.Cells.Item(114, 7).Value = 2
$ws.Cells.Item(114, 8).Value = "Everton"
$ws.Cells.Item(114, 9).Value = 3
$ws.Cells.Item(114, 10).Value = 2.09
$ws.Cells.Item(114, 11).Value = "28/10/2023 20:02"
$ws.Cells.Item(114, 12).Value = 2.71
$ws.Cells.Item(114, 13).Value = "11/11/2023 15:58"
$ws.Cells.Item(114, 14).Value = 3.34
$ws.Cells.Item(114, 15).Value = "28/10/2023 20:02"
$ws.Cells.Item(114, 16).Value = 3.12
$ws.Cells.Item(114, 17).Value = "11/11/2023 15:56"
$ws.Cells.Item(114, 18).Value = 3.8
$ws.Cells.Item(114, 19).Value = "28/10/2023 20:02"
$ws.Cells.Item(114, 20).Value = 2.96
$ws.Cells.Item(114, 21).Value = "11/11/2023 15:59"
$ws.Cells.Item(114, 22).Value = "https://www.betexplorer.com/football/england/premier-league/crystal-palace-everton/rZrW8iVi/"
# --- Append 4 new match rows (117-120), copying A/E cell formatting from the last existing row ---
# Row 117 (new)
$ws.Cells.Item(116, 1).Copy()
$ws.Cells.Item(117, 1).PasteSpecial(-4122)
$ws.Cells.Item(116, 5).Copy()
$ws.Cells.Item(117, 5).PasteSpecial(-4122)
$ws.Cells.Item(117, 1).Value = 116
$ws.Cells.Item(117, 2).Value = "england"
$ws.Cells.Item(117, 3).Value = "premier-league"
$ws.Cells.Item(117, 4).Value = "2023-2024"
$ws.Cells.Item(117, 5).Value = 45242.625
$ws.Cells.Item(117, 6).Value = "Aston Villa"
$ws.Cells.Item(117, 7).Value = 3
$ws.Cells.Item(117, 8).Value = "Fulham"
$ws.Cells.Item(117, 9).Value = 1
$ws.Cells.Item(117, 10).Value = 1.59
$ws.Cells.Item(117, 11).Value = "29/10/2023 00:02"
$ws.Cells.Item(117, 12).Value = 1.64
$ws.Cells.Item(117, 13).Value = "12/11/2023 14:58"
$ws.Cells.Item(117, 14).Value = 4.16
$ws.Cells.Item(117, 15).Value = "29/10/2023 00:02"
$ws.Cells.Item(117, 16).Value = 4.22
$ws.Cells.Item(117, 17).Value = "12/11/2023 14:55"
$ws.Cells.Item(117, 18).Value = 5.06
$ws.Cells.Item(117, 19).Value = "29/10/2023 00:02"
$ws.Cells.Item(117, 20).Value = 5.4
$ws.Cells.Item(117, 21).Value = "12/11/2023 14:58"
$ws.Cells.Item(117, 22).Value = "https://www.betexplorer.com/football/england/premier-league/aston-villa-fulham/EJ4IkUUT/"

# Row 118 (new)
$ws.Cells.Item(116, 1).Copy()
$ws.Cells.Item(118, 1).PasteSpecial(-4122)
$ws.Cells.Item(116, 5).Copy()
$ws.Cells.Item(118, 5).PasteSpecial(-4122)
$ws.Cells.Item(118, 1).Value = 117
$ws.Cells.Item(118, 2).Value = "england"
$ws.Cells.Item(118, 3).Value = "premier-league"
$ws.Cells.Item(118, 4).Value = "2023-2024"
$ws.Cells.Item(118, 5).Value = 45242.625
$ws.Cells.Item(118, 6).Value = "Brighton"
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = "Sheffield Utd"
$ws.Cells.Item(118, 9).Value = 1
$ws.Cells.Item(118, 10).Value = 1.23
$ws.Cells.Item(118, 11).Value = "29/10/2023 11:22"
$ws.Cells.Item(118, 12).Value = 1.27
$ws.Cells.Item(118, 13).Value = "12/11/2023 14:55"
$ws.Cells.Item(118, 14).Value = 6.85
$ws.Cells.Item(118, 15).Value = "29/10/2023 11:22"
$ws.Cells.Item(118, 16).Value = 6.23
$ws.Cells.Item(118, 17).Value = "12/11/2023 14:59"
$ws.Cells.Item(118, 18).Value = 11.63
$ws.Cells.Item(118, 19).Value = "29/10/2023 11:22"
$ws.Cells.Item(118, 20).Value = 11.19
$ws.Cells.Item(118, 21).Value = "12/11/2023 14:59"
$ws.Cells.Item(118, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-sheffield-utd/GSjpolwo/"

# Row 119 (new)
$ws.Cells.Item(116, 1).Copy()
$ws.Cells.Item(119, 1).PasteSpecial(-4122)
$ws.Cells.Item(116, 5).Copy()
$ws.Cells.Item(119, 5).PasteSpecial(-4122)
$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = "england"
$ws.Cells.Item(119, 3).Value = "premier-league"
$ws.Cells.Item(119, 4).Value = "2023-2024"
$ws.Cells.Item(119, 5).Value = 45242.625
$ws.Cells.Item(119, 6).Value = "Liverpool"
$ws.Cells.Item(119, 7).Value = 3
$ws.Cells.Item(119, 8).Value = "Brentford"
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 1.48
$ws.Cells.Item(119, 11).Value = "28/10/2023 22:02"
$ws.Cells.Item(119, 12).Value = 1.4
$ws.Cells.Item(119, 13).Value = "12/11/2023 14:56"
$ws.Cells.Item(119, 14).Value = 4.87
$ws.Cells.Item(119, 15).Value = "28/10/2023 22:02"
$ws.Cells.Item(119, 16).Value = 5.51
$ws.Cells.Item(119, 17).Value = "12/11/2023 14:56"
$ws.Cells.Item(119, 18).Value = 6.16
$ws.Cells.Item(119, 19).Value = "28/10/2023 22:02"
$ws.Cells.Item(119, 20).Value = 7.48
$ws.Cells.Item(119, 21).Value = "12/11/2023 14:59"
$ws.Cells.Item(119, 22).Value = "https://www.betexplorer.com/football/england/premier-league/liverpool-brentford/QcYz8Bpb/"

# Row 120 (new)
$ws.Cells.Item(116, 1).Copy()
$ws.Cells.Item(120, 1).PasteSpecial(-4122)
$ws.Cells.Item(116, 5).Copy()
$ws.Cells.Item(120, 5).PasteSpecial(-4122)
$ws.Cells.Item(120, 1).Value = 119
$ws.Cells.Item(120, 2).Value = "england"
$ws.Cells.Item(120, 3).Value = "premier-league"
$ws.Cells.Item(120, 4).Value = "2023-2024"
$ws.Cells.Item(120, 5).Value = 45242.625
$ws.Cells.Item(120, 6).Value = "West Ham"
$ws.Cells.Item(120, 7).Value = 3
$ws.Cells.Item(120, 8).Value = "Nottingham"
$ws.Cells.Item(120, 9).Value = 2
$ws.Cells.Item(120, 10).Value = 1.82
$ws.Cells.Item(120, 11).Value = "28/10/2023 22:02"
$ws.Cells.Item(120, 12).Value = 1.81
$ws.Cells.Item(120, 13).Value = "12/11/2023 14:35"
$ws.Cells.Item(120, 14).Value = 3.74
$ws.Cells.Item(120, 15).Value = "28/10/2023 22:02"
$ws.Cells.Item(120, 16).Value = 3.84
$ws.Cells.Item(120, 17).Value = "12/11/2023 14:35"
$ws.Cells.Item(120, 18).Value = 4.46
$ws.Cells.Item(120, 19).Value = "28/10/2023 22:02"
$ws.Cells.Item(120, 20).Value = 4.59
$ws.Cells.Item(120, 21).Value = "12/11/2023 14:58"
$ws.Cells.Item(120, 22).Value = "https://www.betexplorer.com/football/england/premier-league/west-ham-nottingham/YNyq6kFA/"
